$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2: clear "To Do" value, keep just formatting (no border fill, applyFill)
$ws.Range("F2").ClearContents()

# F4, F6, F7, F12, F13, F14, F15: change "To Do" -> "Done"
$ws.Range("F4").Value = "Done"
$ws.Range("F6").Value = "Done"
$ws.Range("F7").Value = "Done"
$ws.Range("F12").Value = "Done"
$ws.Range("F13").Value = "Done"
$ws.Range("F14").Value = "Done"
$ws.Range("F15").Value = "Done"

# D10: change "M" -> "S"
$ws.Range("D10").Value = "S"

# Update sheet view: topLeftCell B1, selection G15
$ws.Range("G15").Select()
$excel.ActiveWindow.ScrollColumn = 2
